# Update "想去人数" (want-to-go count) values in column F for rows 2-5
# on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 525
    $ws.Range("F3").Value = 3473
    $ws.Range("F4").Value = 95
    $ws.Range("F5").Value = 679
}
